# correction de l'affichage dans backend
# - "monji" -> "ahmed" on the second reservation row
# - the former "ahmed" row now shows "karoui" again
# - eight new reservation rows (all "karoui") are appended, each with a
#   precise backend timestamp in column B and a "Confirmé" status in
#   column C (highlighted the same way the existing status fill is used)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix existing rows -----------------------------------------------
$ws.Cells.Item(4, 1).Value = "ahmed"
$ws.Cells.Item(5, 1).Value = "karoui"

# --- append the new backend-generated reservations --------------------
$timestamps = @(
    "2025-03-06T12:34:49.917466400",
    "2025-03-06T12:40:07.667345",
    "2025-03-06T12:44:21.003296400",
    "2025-03-06T12:44:34.228014700",
    "2025-03-06T12:53:51.661579400",
    "2025-03-06T12:58:30.918849200",
    "2025-03-06T13:01:30.447345700",
    "2025-03-06T13:02:13.757866900"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = 6 + $i

    $ws.Cells.Item($row, 1).Value = "karoui"
    $ws.Cells.Item($row, 2).Value = $timestamps[$i]
    $ws.Cells.Item($row, 2).HorizontalAlignment = -4108  # xlCenter, matches the other date cells
    $ws.Cells.Item($row, 3).Value = "Confirmé"

    # reuse the same highlighted-status look already used in the
    # workbook's style table (a solid fill, indexed colour 42)
    $ws.Cells.Item($row, 3).Interior.ColorIndex = 42
    $ws.Cells.Item($row, 3).Interior.Pattern = 1
}

# column B now needs to fit the long timestamp strings (closest width the
# host's character-based ColumnWidth metric can express vs. the ~30.84
# "characters" OOXML width recorded by the original backend export)
$ws.Columns.Item(2).ColumnWidth = 30
